# COMPADRE Kendall protocol - add drop-down lists for fixed-response fields.
#
# 1. Rename Sheet1 -> "Data entry"
# 2. Add a new sheet "Data validation" right after it, and move the
#    validation-source lists that used to live in the far-right columns of
#    "Data entry" (M, O, Q, S, U, AA, AI) into compact columns A, C, E, G, I, K
#    on the new sheet, each with its own header + Excel Table.
# 3. Clear the old validation-source cells from "Data entry".
# 4. Point the relevant "Data entry" columns at the new lists via data
#    validation drop-downs.

$wb  = $excel.ActiveWorkbook
$de  = $wb.Worksheets.Item(1)
$de.Name = "Data entry"

$dv = $wb.Worksheets.Add($null, $de)
$dv.Name = "Data validation"

# ---------------------------------------------------------------------
# Populate "Data validation" sheet
# ---------------------------------------------------------------------

# Column A: MatrixModified
$dv.Range("A1").Value = "MatrixModified"
$dv.Range("A2").Value = "Yes"
$dv.Range("A3").Value = "No"
$dv.Range("A4").Value = "No Matrix"

# Column C: MatrixIDfield
$dv.Range("C1").Value = "MatrixIDfield"
$dv.Range("C2").Value = "SpeciesAuthor"
$dv.Range("C3").Value = "MatrixPopulation"
$dv.Range("C4").Value = "MatrixTreatment"
$dv.Range("C5").Value = "MatrixStartYear"
$dv.Range("C6").Value = "MatrixStartSeason"
$dv.Range("C7").Value = "MatrixStartMonth"
$dv.Range("C8").Value = "MatrixEndYear"
$dv.Range("C9").Value = "MatrixEndSeason"
$dv.Range("C10").Value = "MatrixEndMonth"
$dv.Range("C11").Value = "Observation"

# Column E: CensusType
$dv.Range("E1").Value = "CensusType"
$dv.Range("E2").Value = "Pre"
$dv.Range("E3").Value = "Post"
$dv.Range("E4").Value = "Post+"
$dv.Range("E5").Value = "Mid"
$dv.Range("E6").Value = "Flow"
$dv.Range("E7").Value = "Ambiguous"
$dv.Range("E8").Value = "Unknown"

# Column G: SurvInRep
$dv.Range("G1").Value = "SurvInRep"
$dv.Range("G2").Value = "None"
$dv.Range("G3").Value = "Offspring"
$dv.Range("G4").Value = "Parent"
$dv.Range("G5").Value = "OffsetParent"
$dv.Range("G6").Value = "Parent | EarlyOffspring"
$dv.Range("G7").Value = "PartialParent | PartialOffspring"
$dv.Range("G8").Value = "Unknown"

# Column I: ReproWithMaturation
$dv.Range("I1").Value = "ReproWithMaturation"
$dv.Range("I2").Value = "Yes"
$dv.Range("I3").Value = "No"
$dv.Range("I4").Value = "Unknown"

# Column K: GrowthTransition
$dv.Range("K1").Value = "GrowthTransition"
$dv.Range("K2").Value = "1/Tbar"
$dv.Range("K3").Value = "Cohort"
$dv.Range("K4").Value = "Caswell6.103"
$dv.Range("K5").Value = "Observed"
$dv.Range("K6").Value = "Unrolled"
$dv.Range("K7").Value = "Variable"
$dv.Range("K8").Value = "NegativeBinomial"
$dv.Range("K9").Value = "Other"
$dv.Range("K10").Value = "Unknown"

# Header row is bold, including the blank spacer columns (B, D, F, H, J)
$dv.Range("A1:K1").Font.Bold = $true

# Column widths (approximate fit to content, matching manual layout)
$dv.Columns.Item(1).ColumnWidth = 15.5
$dv.Columns.Item(2).ColumnWidth = 6.5
$dv.Columns.Item(3).ColumnWidth = 16.5
$dv.Columns.Item(4).ColumnWidth = 6.83
$dv.Columns.Item(5).ColumnWidth = 12.33
$dv.Columns.Item(6).ColumnWidth = 7.83
$dv.Columns.Item(7).ColumnWidth = 26.5
$dv.Columns.Item(8).ColumnWidth = 6.83
$dv.Columns.Item(9).ColumnWidth = 20.67
$dv.Columns.Item(10).ColumnWidth = 5.33
$dv.Columns.Item(11).ColumnWidth = 17

# ---------------------------------------------------------------------
# Turn each list into an Excel Table (so the named dropdown sources are
# robust / self-expanding, matching the authored workbook).
# ---------------------------------------------------------------------
$t1 = $dv.ListObjects.Add(1, $dv.Range("A1:A4"), $null, 1)
$t2 = $dv.ListObjects.Add(1, $dv.Range("C1:C11"), $null, 1)
$t3 = $dv.ListObjects.Add(1, $dv.Range("E1:E8"), $null, 1)
$t4 = $dv.ListObjects.Add(1, $dv.Range("G1:G8"), $null, 1)
$t5 = $dv.ListObjects.Add(1, $dv.Range("I1:I4"), $null, 1)
$t6 = $dv.ListObjects.Add(1, $dv.Range("K1:K10"), $null, 1)

# Rename from the end backwards to avoid clobbering earlier default names.
$t6.Name = "Table7"
$t5.Name = "Table6"
$t4.Name = "Table5"
$t3.Name = "Table4"
# $t2 keeps its default name "Table2"
# $t1 keeps its default name "Table1"

# ---------------------------------------------------------------------
# Clear the old validation-source data from "Data entry"
# ---------------------------------------------------------------------
$de.Range("M2:AJ11").ClearContents()

# ---------------------------------------------------------------------
# Wire up drop-down data validation on "Data entry", pointing at the new
# "Data validation" sheet ranges.
# ---------------------------------------------------------------------
$de.Range("M2:M201").Validation.Add(3, 1, 1, "='Data validation'!`$A`$2:`$A`$4")
$de.Range("O2:O201").Validation.Add(3, 1, 1, "='Data validation'!`$C`$2:`$C`$11")
$de.Range("Q2:Q201").Validation.Add(3, 1, 1, "='Data validation'!`$C`$2:`$C`$11")
$de.Range("S2:S201").Validation.Add(3, 1, 1, "='Data validation'!`$E`$2:`$E`$8")
$de.Range("U2:U201").Validation.Add(3, 1, 1, "='Data validation'!`$G`$2:`$G`$8")
$de.Range("AA2:AA201").Validation.Add(3, 1, 1, "='Data validation'!`$I`$2:`$I`$4")
$de.Range("AI2:AI201").Validation.Add(3, 1, 1, "='Data validation'!`$K`$2:`$K`$10")

foreach ($rng in @("M2:M201","O2:O201","Q2:Q201","S2:S201","U2:U201","AA2:AA201","AI2:AI201")) {
  $v = $de.Range($rng).Validation
  $v.IgnoreBlank = $true
  $v.InCellDropdown = $true
  $v.ShowInput = $true
  $v.ShowError = $true
}

# ---------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------
$dv.Range("C17").Select()
$de.Activate()
$de.Range("A2").Select()
